$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 199, shifting existing rows 199-204 down to 200-205.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new weekly price record.
$ws.Cells.Item(199, 1).Value = 9
$ws.Cells.Item(199, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(199, 3).Value = "Metropolitana"
$ws.Cells.Item(199, 4).Value = 44448
$ws.Cells.Item(199, 5).Value = 13
$ws.Cells.Item(199, 6).Value = 100112044
$ws.Cells.Item(199, 7).Value = "Perejil"
$ws.Cells.Item(199, 8).Value = "Sin especificar"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 142
$ws.Cells.Item(199, 11).Value = 9000
$ws.Cells.Item(199, 12).Value = 10000
$ws.Cells.Item(199, 13).Value = 9500
$ws.Cells.Item(199, 14).Value = "`$/docena de atados"
$ws.Cells.Item(199, 15).Value = "Región Metropolitana"
$ws.Cells.Item(199, 16).Value = 3167
$ws.Cells.Item(199, 17).Value = 3
$ws.Cells.Item(199, 18).Value = "Hortaliza"
